$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171132445335388
$ws.Range("B1").Value = 2.438088655471802
$ws.Range("D1").Value = 2.362921476364136
$ws.Range("E1").Value = 1.238832354545593
